$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.md"
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("A3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.md"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.md"
$zhcn.Range("C2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.5d1d59fb21b8c08a1c21d2083a422304aba18e35.zh-cn.xlf"
$zhcn.Range("D2").Value = "2016-01-18 02:04:35"
$zhcn.Range("E2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.md"
$zhcn.Range("F2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.5d1d59fb21b8c08a1c21d2083a422304aba18e35.zh-cn.xlf"
$zhcn.Range("G2").Value = "2016-01-18 02:05:23"
$zhcn.Range("H2").Value = "Include"

$zhcn.Range("A3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.md"
$zhcn.Range("C3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.df3a4e6a6af4098ef13c7aa20cb724752ea738b3.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-18 02:02:45"
$zhcn.Range("E3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.md"
$zhcn.Range("F3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.df3a4e6a6af4098ef13c7aa20cb724752ea738b3.zh-cn.xlf"
$zhcn.Range("G3").Value = "2016-01-18 02:03:31"
$zhcn.Range("H3").Value = "Include"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.md"
$dede.Range("C2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.5d1d59fb21b8c08a1c21d2083a422304aba18e35.de-de.xlf"
$dede.Range("D2").Value = "2016-01-18 02:04:47"
$dede.Range("E2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.md"
$dede.Range("F2").Value = "24b49878-1168-41df-9209-d28d7c304a0f.5d1d59fb21b8c08a1c21d2083a422304aba18e35.de-de.xlf"
$dede.Range("G2").Value = "2016-01-18 02:05:45"
$dede.Range("H2").Value = "Include"

$dede.Range("A3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.md"
$dede.Range("C3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.df3a4e6a6af4098ef13c7aa20cb724752ea738b3.de-de.xlf"
$dede.Range("D3").Value = "2016-01-18 02:02:58"
$dede.Range("E3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.md"
$dede.Range("F3").Value = "f1119a0a-5a0e-4964-908d-19df625c50a4.df3a4e6a6af4098ef13c7aa20cb724752ea738b3.de-de.xlf"
$dede.Range("G3").Value = "2016-01-18 02:03:52"
$dede.Range("H3").Value = "Include"
